# Update the Sprint Backlog / Burndown sheet:
#  - Week 1 ("D") amounts for the first two tasks are now tracked as 0
#    (explicit zero, rather than a blank cell).
#  - Move/restore the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("B3").Select()
